$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 419, pushing existing rows 419:489 down to 420:490
$ws.Rows("419:419").Insert()

# Populate the newly inserted row 419 with the new weekly data record
$ws.Range("A419").Value = 11
$ws.Range("B419").Value = "Vega Monumental Concepción"
$ws.Range("C419").Value = "Bíobío"
$ws.Range("D419").Value = 45258
$ws.Range("E419").Value = 8
$ws.Range("F419").Value = 100112009
$ws.Range("G419").Value = "Acelga"
$ws.Range("H419").Value = "Sin especificar"
$ws.Range("I419").Value = "Primera"
$ws.Range("J419").Value = 450
$ws.Range("K419").Value = 600
$ws.Range("L419").Value = 650
$ws.Range("M419").Value = 617
$ws.Range("N419").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O419").Value = "Región de Ñuble"
$ws.Range("P419").Value = 617
$ws.Range("Q419").Value = 1
$ws.Range("R419").Value = "Hortaliza"
